# Add 2022-Q3 data
# -------------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" right after the "总计" sheet,
#    by duplicating the existing "2022-Q2" sheet (so it inherits the exact
#    same layout/formatting) and then overwriting its data.
# 2) Insert a new row into "总计" (right after its header row) summarizing
#    the 2022-Q3 totals, shifting the existing rows down, copying the
#    formatting from the row below so the new row matches the rest of the
#    table, and renumbering the running index in column A.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---- 1) Duplicate the "2022-Q2" sheet right after "总计" and rename it --
$q2Sheet.Copy($null, $totalSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---- overwrite its data with the 2022-Q3 numbers --------------------------
# columns: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @(0, "270021", "广发聚瑞混合A",     "17.58", "93.91", "4.39", "0.7718", 8),
    @(1, "011136", "广发盛兴混合A",     "16.84", "94.77", "3.63", "0.6113", 9),
    @(2, "010161", "广发瑞安精选股票A", "6.12",  "94.02", "4.13", "0.2528", 8),
    @(3, "011137", "广发盛兴混合C",     "1.74",  "94.77", "3.63", "0.0632", 9),
    @(4, "010026", "广发聚瑞混合C",     "0.50",  "93.91", "4.39", "0.0220", 8),
    @(5, "010162", "广发瑞安精选股票C", "0.53",  "94.02", "4.13", "0.0219", 8)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]

    $q3.Cells.Item($r, 1).Value = $row[0]

    # Columns 基金代码/基金规模/股票总仓位/仓位占比/持有市值 are plain text
    # (leading/trailing zeros such as "010026" / "0.50" / "0.0220" must
    # survive), so force text interpretation while assigning ...
    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[6]

    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# ... then restore the original (default) cell formatting from the
# still-untouched "2022-Q2" sheet so the new sheet's look matches the rest
# of the workbook exactly.
$q2Sheet.Range("B2:G7").Copy()
$q3.Range("B2:G7").PasteSpecial($xlPasteFormats)

# ---- 2) Insert the 2022-Q3 summary row into "总计" ------------------------
$totalSheet.Rows.Item(2).Insert()

# match the formatting of the data rows (border-less, not bold) instead of
# whatever got pulled in by the row insert
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial($xlPasteFormats)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 1.74

# renumber the running index in column A (0,1,2,...) for every data row
# now that an extra row was inserted at the top
$lastRow = $totalSheet.Cells.Item(1, 2).End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
